# Apply load forecast update: shift dates from day 45431 to 45435-45436,
# update load values, and extend data through row 98 (new rows 94-98 added).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$aValues = @(45435,45435.01041666666,45435.02083333334,45435.03125,45435.04166666666,45435.05208333334,45435.0625,45435.07291666666,45435.08333333334,45435.09375,45435.10416666666,45435.11458333334,45435.125,45435.13541666666,45435.14583333334,45435.15625,45435.16666666666,45435.17708333334,45435.1875,45435.19791666666,45435.20833333334,45435.21875,45435.22916666666,45435.23958333334,45435.25,45435.26041666666,45435.27083333334,45435.28125,45435.29166666666,45435.30208333334,45435.3125,45435.32291666666,45435.33333333334,45435.34375,45435.35416666666,45435.36458333334,45435.375,45435.38541666666,45435.39583333334,45435.40625,45435.41666666666,45435.42708333334,45435.4375,45435.44791666666,45435.45833333334,45435.46875,45435.47916666666,45435.48958333334,45435.5,45435.51041666666,45435.52083333334,45435.53125,45435.54166666666,45435.55208333334,45435.5625,45435.57291666666,45435.58333333334,45435.59375,45435.60416666666,45435.61458333334,45435.625,45435.63541666666,45435.64583333334,45435.65625,45435.66666666666,45435.67708333334,45435.6875,45435.69791666666,45435.70833333334,45435.71875,45435.72916666666,45435.73958333334,45435.75,45435.76041666666,45435.77083333334,45435.78125,45435.79166666666,45435.80208333334,45435.8125,45435.82291666666,45435.83333333334,45435.84375,45435.85416666666,45435.86458333334,45435.875,45435.88541666666,45435.89583333334,45435.90625,45435.91666666666,45435.92708333334,45435.9375,45435.94791666666,45435.95833333334,45435.96875,45435.97916666666,45435.98958333334,45436)
$bValues = @(5330,5290,5250,5210,5180,5150,5130,5120,5120,5130,5130,5140,5140,5140,5140,5150,5170,5210,5260,5330,5420,5530,5660,5800,5950,6090,6220,6340,6440,6520,6570,6590,6580,6550,6500,6440,6380,6310,6250,6180,6130,6080,6030,5990,5950,5910,5870,5840,5800,5780,5760,5750,5740,5740,5740,5750,5750,5750,5760,5780,5810,5850,5900,5960,6010,6060,6100,6130,6160,6200,6260,6350,6450,6550,6630,6690,6760,6840,6920,7030,7100,7100,7090,7000,6830,6720,6580,6420,6190,6020,5930,5820,5570,5530,5470,5400,5350)

$startRow = 2
$endRow = $startRow + $aValues.Length - 1

# Build 2D arrays for bulk write
$n = $aValues.Length
$dataArray = New-Object 'object[,]' $n,2
for ($i = 0; $i -lt $n; $i++) {
    $dataArray[$i,0] = $aValues[$i]
    $dataArray[$i,1] = $bValues[$i]
}

$rangeA1 = $ws.Cells.Item($startRow, 1)
$rangeB2 = $ws.Cells.Item($endRow, 2)
$targetRange = $ws.Range($rangeA1, $rangeB2)
$targetRange.Value = $dataArray

# Make sure the timestamp column keeps (and the newly appended rows gain)
# the same date/time display format used by the existing data (row 93 and
# above), since new rows added past the old A1:B93 extent otherwise default
# to the General format.
$timestampColumn = $ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($endRow, 1))
$timestampColumn.NumberFormat = "YYYY-MM-DD HH:MM:SS"
